# Locate the "Author" styled paragraph that contains "Edison Achalma"
# (the one directly under the document title) and insert a new "Author"
# styled paragraph right after it containing the author's affiliation.

$d = $word.ActiveDocument

$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Range.Style.NameLocal
    $text = $p.Range.Text.Trim()
    if ($styleName -eq "Author" -and $text -eq "Edison Achalma") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Edison Achalma' Author paragraph"
}

# The paragraph immediately following "Edison Achalma" is where we want
# the new paragraph to land. Inserting a paragraph mark after that
# following paragraph's range actually creates a new, empty paragraph
# *between* the two (i.e. right after "Edison Achalma"), leaving the
# "Edison Achalma" paragraph itself untouched.
$followingPara = $d.Paragraphs($targetIndex + 1)
$followingPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($targetIndex + 1)
$newPara.Range.Style = "Author"
$newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
